$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Break up the existing "Application #2" merge (A3:A4) so each row
#    can receive its own, independent value.
# ------------------------------------------------------------------
$ws.Range("A3:A4").UnMerge()

# ------------------------------------------------------------------
# 2. Update the "Versions" row (was row 2 -> ends up row 4) and the
#    "Application" row 1 of 2 (was row 3 -> ends up row 5) in place.
#    Row 2/3 become the brand-new "Non trouve" block.
# ------------------------------------------------------------------
$ws.Range("A2").Value = "Non trouvé`n#1"
$ws.Range("B2").Value = "#1000"
$ws.Range("C2").Value = "Erreur inconnue."

$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "#1001"
$ws.Range("C3").Value = "Fichier non trouvé."

$ws.Range("A4").Value = "Versions`n#2"
$ws.Range("B4").Value = "#2001"
$ws.Range("C4").Value = "Historique des versions introuvable."

# ------------------------------------------------------------------
# 3. Re-merge A2:A3 for the new "Non trouve" block.
# ------------------------------------------------------------------
$ws.Range("A2:A3").Merge()

# ------------------------------------------------------------------
# 4. Insert three fresh rows after row 4 to host the (now 3-row)
#    "Application #3" block, then fill them in and re-merge A5:A7.
# ------------------------------------------------------------------
$ws.Rows.Item("5:7").Insert()

$ws.Range("A5").Value = "Application`n#3"
$ws.Range("B5").Value = "#3001"
$ws.Range("C5").Value = "Application inconnue."

$ws.Range("A6").Value = ""
$ws.Range("B6").Value = "#3002"
$ws.Range("C6").Value = "Contenu de l'application inacessible."

$ws.Range("A7").Value = ""
$ws.Range("B7").Value = "#3003"
$ws.Range("C7").Value = "Impossible de récupérer la correction."

$ws.Range("A5:A7").Merge()

Write-Output "done"
